# Applies the "variable name harmonization with label types" update to the
# var_lab_dict worksheet:
#   - add a "justice$" alias to the "justice-system" row (column D)
#   - insert a new row ("court-of-justice" / "court_of_justice") after
#     "court-of-auditors"
#   - add an "eu_court_of_justice" alias to the "eu-court-of-justice" row
#     (column D)
#   - add a "rely_on" alias to the "trust" row (column C)
#   - insert a new row ("problem" / "prbl") before the row with "services_"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) "justice-system" row (row 57 in the original layout): add a fourth
#    synonym "justice$" in column D.
$ws.Range("D57").Value = "justice$"

# 2) Insert new row "court-of-justice" / "court_of_justice" right after the
#    "court-of-auditors" row (row 106 in the original layout).
$ws.Rows.Item(107).Insert()
$ws.Range("A107").Value = "court-of-justice"
$ws.Range("B107").Value = "court_of_justice"

# 3) "eu-court-of-justice" row (originally row 122, now row 123 after the
#    insert above): add a fourth synonym "eu_court_of_justice" in column D.
$ws.Range("D123").Value = "eu_court_of_justice"

# 4) "trust" row (row 56, unaffected by the insert above): add a third
#    synonym "rely_on" in column C.
$ws.Range("C56").Value = "rely_on"

# 5) Insert new row "problem" / "prbl" above the row holding "services_"
#    (row 44 in the original layout).
$ws.Rows.Item(44).Insert()
$ws.Range("B44").Value = "prbl"
$ws.Range("A44").Value = "problem"
